$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (new): TC003 "Checkbox State" ---
$ws.Range("A4").Value = "TC003"
$ws.Range("B4").Value = "Checkbox State"
$ws.Range("C4").Value = "Verify checkbox state changes correctly after pressing checking/unchecking."
$ws.Range("D4").Value = "Navigate to the Check Box page"

# --- Row 3 (TC002): trim the last test step out of the "Test Steps" cell,
#     shrinking the wrapped text from 4 lines to 3 lines. Row 4 reuses the
#     same (now shorter) steps text. The wrap style (s=1) is already applied
#     to E3, so only the text itself needs to change there. ---
$ws.Range("E3").Value = "1. Expand the Home Folder.`n2. Click the checkbox of ""Home"" checkbox.`n3. Click the checkbox of ""Documents"" checkbox."
$ws.Range("E4").Value = "1. Expand the Home Folder.`n2. Click the checkbox of ""Home"" checkbox.`n3. Click the checkbox of ""Documents"" checkbox."
$ws.Range("E4").WrapText = $true

$ws.Range("F4").Value = "Checked state of checkbox should be displayed correctly."

$ws.Range("F3").Value = "Selected/deselected items should display in the output list."

$ws.Range("G4").Value = "Pass"
$ws.Range("G4").HorizontalAlignment = -4108
$ws.Range("G4").VerticalAlignment = -4108

# Row heights shrink/match the (now 3-line) wrapped Test Steps text.
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 75

# --- Selection moves down to D12 ---
$ws.Range("D12").Select()
